$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 13 data: plant and terra cota holder from Lowes
$ws.Range("A13").Value = "Plant and terra cota holder"
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = 15
$ws.Range("D13").Value = "Lowes"

# Update the total formula to include the new row
$ws.Range("C15").Formula = "=SUM(C2:C13)"

# Update selection to A13
$ws.Range("A13").Select()
